$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1331.3334
$ws.Range("I5").Value = 997
$ws.Range("K5").Value = 997
$ws.Range("M5").Value = -882
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4825
$ws.Range("H70").Value = 1476.7333
$ws.Range("J70").Value = 1244.1666
$ws.Range("L70").Value = 3732.4998
$ws.Range("N70").Value = -4272.4998
$ws.Range("H73").Value = 1476.7333
$ws.Range("J73").Value = 1244.1666
$ws.Range("L73").Value = 3732.4998
$ws.Range("N73").Value = -5604.4998
$ws.Range("H116").Value = 16350
$ws.Range("I116").Value = 4900
$ws.Range("K116").Value = 4900
$ws.Range("M116").Value = -1458
$ws.Range("H137").Value = 1401.6818
$ws.Range("I137").Value = 1299.3125
$ws.Range("J137").Value = 1674.6666
$ws.Range("K137").Value = 3897.9375
$ws.Range("L137").Value = 5023.9998
$ws.Range("M137").Value = -1347.9375
$ws.Range("N137").Value = -10123.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 29413900
$ws.Range("I2").Value = 33335420
$ws.Range("K2").Value = 33335420
$ws.Range("M2").Value = -33335307
$ws.Range("H32").Value = 3676.8333
$ws.Range("I32").Value = 2318.4443
$ws.Range("K32").Value = 2318.4443
$ws.Range("M32").Value = -2031.4443
$ws.Range("H45").Value = 11646.071
$ws.Range("J45").Value = 3131.111
$ws.Range("L45").Value = 3131.111
$ws.Range("N45").Value = -3885.111
$ws.Range("H55").Value = 24999
$ws.Range("J55").Value = 24999
$ws.Range("L55").Value = 24999
$ws.Range("N55").Value = -25629
$ws.Range("H61").Value = 6186.875
$ws.Range("I61").Value = 4349.25
$ws.Range("J61").Value = 8024.5
$ws.Range("K61").Value = 4349.25
$ws.Range("L61").Value = 8024.5
$ws.Range("M61").Value = -4137.25
$ws.Range("N61").Value = -8448.5
$ws.Range("H110").Value = 1458.317
$ws.Range("I110").Value = 1522.8158
$ws.Range("J110").Value = 641.3333
$ws.Range("K110").Value = 1522.8158
$ws.Range("L110").Value = 641.3333
$ws.Range("M110").Value = 522.1841999999999
$ws.Range("N110").Value = -4731.3333
$ws.Range("H116").Value = 29413900
$ws.Range("I116").Value = 33335420
$ws.Range("K116").Value = 33335420
$ws.Range("M116").Value = -33333126
$ws.Range("H136").Value = 6186.875
$ws.Range("I136").Value = 4349.25
$ws.Range("J136").Value = 8024.5
$ws.Range("K136").Value = 13047.75
$ws.Range("L136").Value = 24073.5
$ws.Range("M136").Value = -10497.75
$ws.Range("N136").Value = -29173.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 29413900
$ws.Range("I3").Value = 33335420
$ws.Range("K3").Value = 33335420
$ws.Range("M3").Value = -33335306
$ws.Range("H58").Value = 16750
$ws.Range("J58").Value = 16750
$ws.Range("L58").Value = 16750
$ws.Range("N58").Value = -17338
$ws.Range("H59").Value = 224833.33
$ws.Range("I59").Value = 100000
$ws.Range("J59").Value = 287250
$ws.Range("K59").Value = 100000
$ws.Range("L59").Value = 287250
$ws.Range("M59").Value = -99153
$ws.Range("N59").Value = -288944

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4618.778
$ws.Range("I6").Value = 8114
$ws.Range("K6").Value = 8114
$ws.Range("M6").Value = -8001
$ws.Range("H16").Value = 2759
$ws.Range("I16").Value = 3964.6667
$ws.Range("J16").Value = 2357.111
$ws.Range("K16").Value = 3964.6667
$ws.Range("L16").Value = 2357.111
$ws.Range("M16").Value = -3677.6667
$ws.Range("N16").Value = -2931.111
$ws.Range("H43").Value = 84599.875
$ws.Range("J43").Value = 84599.875
$ws.Range("L43").Value = 84599.875
$ws.Range("N43").Value = -84967.875
$ws.Range("H62").Value = 333338620
$ws.Range("I62").Value = 500003940
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 500003940
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -500003316
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 333338620
$ws.Range("I65").Value = 500003940
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 2500019700
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -2500016580
$ws.Range("N65").Value = -46240
$ws.Range("H101").Value = 84599.875
$ws.Range("J101").Value = 84599.875
$ws.Range("L101").Value = 84599.875
$ws.Range("N101").Value = -91089.875
$ws.Range("H113").Value = 2759
$ws.Range("I113").Value = 3964.6667
$ws.Range("J113").Value = 2357.111
$ws.Range("K113").Value = 3964.6667
$ws.Range("L113").Value = 2357.111
$ws.Range("M113").Value = -1794.6667
$ws.Range("N113").Value = -6697.111

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 296.2
$ws.Range("I103").Value = 177.22223
$ws.Range("K103").Value = 531.66669
$ws.Range("M103").Value = 347.33331

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 48999.5
$ws.Range("J95").Value = 48999.5
$ws.Range("L95").Value = 48999.5
$ws.Range("N95").Value = -54491.5
$ws.Range("H122").Value = 3103.147
$ws.Range("I122").Value = 3261.353
$ws.Range("J122").Value = 2944.9412
$ws.Range("K122").Value = 9784.059000000001
$ws.Range("L122").Value = 8834.8236
$ws.Range("M122").Value = -7334.059000000001
$ws.Range("N122").Value = -13734.8236

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4202.4165
$ws.Range("I7").Value = 4140.091
$ws.Range("J7").Value = 4888
$ws.Range("K7").Value = 4140.091
$ws.Range("L7").Value = 4888
$ws.Range("M7").Value = -4028.091
$ws.Range("N7").Value = -5112
$ws.Range("H40").Value = 15893.208
$ws.Range("I40").Value = 17372.47
$ws.Range("J40").Value = 12300.714
$ws.Range("K40").Value = 17372.47
$ws.Range("L40").Value = 12300.714
$ws.Range("M40").Value = -17236.47
$ws.Range("N40").Value = -12572.714
$ws.Range("H61").Value = 8027.643
$ws.Range("I61").Value = 10810
$ws.Range("J61").Value = 1071.75
$ws.Range("K61").Value = 10810
$ws.Range("L61").Value = 1071.75
$ws.Range("M61").Value = -10608
$ws.Range("N61").Value = -1475.75
$ws.Range("H113").Value = 8027.643
$ws.Range("I113").Value = 10810
$ws.Range("J113").Value = 1071.75
$ws.Range("K113").Value = 10810
$ws.Range("L113").Value = 1071.75
$ws.Range("M113").Value = -8640
$ws.Range("N113").Value = -5411.75
$ws.Range("H122").Value = 8064.1562
$ws.Range("I122").Value = 8102.3477
$ws.Range("J122").Value = 7966.5557
$ws.Range("K122").Value = 24307.0431
$ws.Range("L122").Value = 23899.6671
$ws.Range("M122").Value = -21857.0431
$ws.Range("N122").Value = -28799.6671
$ws.Range("H126").Value = 4202.4165
$ws.Range("I126").Value = 4140.091
$ws.Range("J126").Value = 4888
$ws.Range("K126").Value = 12420.273
$ws.Range("L126").Value = 14664
$ws.Range("M126").Value = -9950.273000000001
$ws.Range("N126").Value = -19604

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 249
$ws.Range("I55").Value = 48
$ws.Range("J55").Value = 450
$ws.Range("K55").Value = 48
$ws.Range("L55").Value = 450
$ws.Range("M55").Value = 229
$ws.Range("N55").Value = -1004
$ws.Range("H107").Value = 27778324
$ws.Range("I107").Value = 638
$ws.Range("J107").Value = 55556012
$ws.Range("K107").Value = 1914
$ws.Range("L107").Value = 166668036
$ws.Range("M107").Value = 6
$ws.Range("N107").Value = -166671876
$ws.Range("H113").Value = 541.8261
$ws.Range("I113").Value = 353.2
$ws.Range("J113").Value = 895.5
$ws.Range("K113").Value = 1059.6
$ws.Range("L113").Value = 2686.5
$ws.Range("M113").Value = 1110.4
$ws.Range("N113").Value = -7026.5
$ws.Range("H126").Value = 10473.2
$ws.Range("I126").Value = 9191.909
$ws.Range("K126").Value = 27575.727
$ws.Range("M126").Value = -25105.727
